$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.105.97"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.758.30"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'576.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").Value = "'159.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("E9").Value = "  -4.23%  "
$ws.Range("D10").Value = "'5.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.05%  "
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("D13").Value = "3.247.15"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "63.729.86"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("E16").Value = "  -5.54%  "
$ws.Range("D17").Value = "2.763.82"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("E19").Value = "  -5.46%  "
$ws.Range("D20").Value = "'359.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "'6.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.28%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  -8.54%  "
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("E25").Value = "  -4.10%  "
$ws.Range("E26").Value = "  -3.98%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "0.0₃0905"
$ws.Range("E28").Value = "  -7.00%  "
$ws.Range("D29").Value = "'7.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.71%  "
$ws.Range("D32").Value = "'170.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").Value = "'4.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.40%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "'347.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").Value = "'6.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "'4.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("D42").Value = "'39.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("D43").Value = "'21.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.54%  "
$ws.Range("D44").Value = "'21.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("D45").Value = "'0.0589"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.77%  "
$ws.Range("D46").Value = "'137.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "'0.0253"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").Value = "'0.628"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  +0.16%  "
